# Actualización automática 2025-06-09 17:00:11
# Registra las nuevas ventas de junio para:
#   - GUERRERO FAREZ FABIAN MAURICIO / BUSTAMANTE ROSERO MARCO TULIO -> 240X120 PORCELANATO : 1026.43
#   - GUERRERO FAREZ FABIAN MAURICIO / PEÑALOZA LOPEZ RONALD STALIN  -> PIEDRA SINTERIZADA   : 997.92
# y actualiza en cascada los totales/porcentajes de cumplimiento dependientes.

$wb = $excel.ActiveWorkbook

$wsGrupo  = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMes    = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumpl  = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# ---------------------------------------------------------------------------
# Hoja "VENTAS POR GRUPO": venta de la categoria 240X120 PORCELANATO (col C)
# para el cliente BUSTAMANTE ROSERO MARCO TULIO (fila 14).
# ---------------------------------------------------------------------------
$wsGrupo.Range("C14").Value2 = 1026.43

# venta de la categoria PIEDRA SINTERIZADA (col K) para el cliente
# PEÑALOZA LOPEZ RONALD STALIN (fila 40).
$wsGrupo.Range("K40").Value2 = 997.92

# Conteo de clientes con venta registrada por categoria (fila 53, "X de 51").
$wsGrupo.Range("C53").Value2 = "3 de 51"
$wsGrupo.Range("K53").Value2 = "2 de 51"

# ---------------------------------------------------------------------------
# Hoja "VENTA MENSUAL": venta de junio (col F) por cliente.
# ---------------------------------------------------------------------------
$wsMes.Range("F14").Value2 = 1000.77
$wsMes.Range("F40").Value2 = 997.92
$wsMes.Range("F53").Value2 = 12847.45

# ---------------------------------------------------------------------------
# Hoja "CUMPLIMIENTO MENSUAL": venta, saldo por cumplir y % de cumplimiento
# por grupo de productos.
# ---------------------------------------------------------------------------

# Fila 2: 240X120 PORCELANATO
$presupuesto2 = $wsCumpl.Range("C2").Value2
$wsCumpl.Range("D2").Value2 = 2042.49
$wsCumpl.Range("E2").Value2 = $presupuesto2 - 2042.49
$wsCumpl.Range("F2").Value2 = 2042.49 / $presupuesto2

# Fila 15: PIEDRA SINTERIZADA
$presupuesto15 = $wsCumpl.Range("C15").Value2
$wsCumpl.Range("D15").Value2 = 1853.28
$wsCumpl.Range("E15").Value2 = $presupuesto15 - 1853.28
$wsCumpl.Range("F15").Value2 = 1853.28 / $presupuesto15

# Fila 19: TOTAL
$presupuesto19 = $wsCumpl.Range("C19").Value2
$wsCumpl.Range("D19").Value2 = 12847.45
$wsCumpl.Range("E19").Value2 = $presupuesto19 - 12847.45
$wsCumpl.Range("F19").Value2 = 12847.45 / $presupuesto19
